# Weekly update: a new price record for the week of 2023-04-17 is inserted
# at row 32 (the data block is sorted most-recent-first), pushing every
# existing record from row 32 down through row 59 down by one row (to
# rows 33-60). The sheet's used range grows from A1:R59 to A1:R60.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 32, shifting rows 32:59 down to 33:60
# (this also pushes the row's formatting down, same as Excel's UI "Insert").
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new weekly record.
$ws.Range("A32").Value = 1
$ws.Range("B32").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C32").Value = "Arica y Parinacota"
$ws.Range("D32").Value = 45033
$ws.Range("E32").Value = 15
$ws.Range("F32").Value = 100112052
$ws.Range("G32").Value = "Albahaca"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 600
$ws.Range("K32").Value = 400
$ws.Range("L32").Value = 500
$ws.Range("M32").Value = 442
$ws.Range("N32").Value = "$/paquete"
$ws.Range("O32").Value = "Región de Arica y Parinacota"
$ws.Range("P32").Value = 442
$ws.Range("Q32").Value = 1
$ws.Range("R32").Value = "Hortaliza"
